$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.1111111111111111
$ws.Range("B1").Value = 0.3333333333333333
$ws.Range("C1").Value = 0.01785714285714286
$ws.Range("D1").Value = 0.01354488130552715
$ws.Range("E1").ClearContents()
$ws.Range("F1").Value = 1
$ws.Range("G1").Value = 5
$ws.Range("H1").Value = 0.2027089762611055
